$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated calcTime (column E) measurements after changing the critic
# network and reducing to 8 actions.
$ws.Range("E3").Value = 0.1479491
$ws.Range("E4").Value = 0.0161172
$ws.Range("E5").Value = 0.023758

$ws.Range("E11").Value = 0.0114814
$ws.Range("E12").Value = 0.0232328
$ws.Range("E13").Value = 0.0299603

$ws.Range("E19").Value = 0.011464
$ws.Range("E20").Value = 0.0118885
$ws.Range("E21").Value = 0.0248873

$ws.Range("E27").Value = 0.0099711
$ws.Range("E28").Value = 0.0109009
$ws.Range("E29").Value = 0.021766
